$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(2).ColumnWidth = 13.22
Write-Host "colwidth set ok"
$ws.Rows.Item(63).Insert()
Write-Host "insert ok"
$ws.Range("A63").Value = 54
Write-Host $ws.Range("A64").Value
